# Refresh the cryptos table with updated prices and 1h volume-change figures.
# Also swaps the Toncoin/Monero entries (rows 29-30) per the upstream data pull.
#
# Numeric-looking price strings (e.g. "299.60") are written with a leading
# apostrophe so Excel keeps them as literal text (preserving trailing zeros)
# instead of silently coercing them to numbers, matching how the sheet
# originally stored these values as text.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.910.13"
$ws.Range("E2").Value = "  +0.59%  "

$ws.Range("D3").Value = "2.295.39"
$ws.Range("E3").Value = "  +0.45%  "

$ws.Range("D5").Value = "'299.60"
$ws.Range("E5").Value = "  -0.60%  "

$ws.Range("D6").Value = "'97.17"
$ws.Range("E6").Value = "  -0.59%  "

$ws.Range("D7").Value = "'0.505"
$ws.Range("E7").Value = "  +0.58%  "

$ws.Range("E8").Value = "  -0.01%  "

$ws.Range("E9").Value = "  +0.23%  "

$ws.Range("D10").Value = "'33.72"
$ws.Range("E10").Value = "  +0.20%  "

$ws.Range("E11").Value = "  +0.36%  "

$ws.Range("D12").Value = "'49.12"
$ws.Range("E12").Value = "  -2.89%  "

$ws.Range("E13").Value = "  +2.96%  "

$ws.Range("D14").Value = "'17.11"
$ws.Range("E14").Value = "  +11.91%  "

$ws.Range("E15").Value = "  +1.68%  "

$ws.Range("D16").Value = "2.646.12"
$ws.Range("E16").Value = "  +0.28%  "

$ws.Range("D17").Value = "2.299.47"
$ws.Range("E17").Value = "  +0.58%  "

$ws.Range("D19").Value = "42.857.09"
$ws.Range("E19").Value = "  +0.68%  "

$ws.Range("D20").Value = "'11.65"
$ws.Range("E20").Value = "  +1.09%  "

$ws.Range("E21").Value = "  +0.55%  "

$ws.Range("D22").Value = "'6.05"
$ws.Range("E22").Value = "  +0.66%  "

$ws.Range("D23").Value = "'67.49"
$ws.Range("E23").Value = "  +1.15%  "

$ws.Range("D24").Value = "'236.54"
$ws.Range("E24").Value = "  +0.63%  "

$ws.Range("E25").Value = "  +5.23%  "

$ws.Range("E26").Value = "  -0.02%  "

$ws.Range("E27").Value = "  -1.58%  "

$ws.Range("B29").Value = "Monero"
$ws.Range("C29").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D29").Value = "'166.32"
$ws.Range("E29").Value = "  +0.88%  "

$ws.Range("B30").Value = "Toncoin"
$ws.Range("C30").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D30").Value = "'2.06"
$ws.Range("E30").Value = "  -9.10%  "

$ws.Range("D31").Value = "'33.73"
$ws.Range("E31").Value = "  -0.08%  "

$ws.Range("D32").Value = "'9.10"

$ws.Range("E33").Value = "  +0.06%  "

$ws.Range("D34").Value = "'4.96"
$ws.Range("E34").Value = "  -0.41%  "

$ws.Range("E35").Value = "  +5.19%  "

$ws.Range("E36").Value = "  +0.69%  "

$ws.Range("E37").Value = "  +3.13%  "

$ws.Range("D38").Value = "'0.0697"
$ws.Range("E38").Value = "  +0.42%  "

$ws.Range("D39").Value = "'2.82"
$ws.Range("E39").Value = "  +0.00%  "

$ws.Range("E40").Value = "  +0.83%  "

$ws.Range("E41").Value = "  -0.35%  "

$ws.Range("E42").Value = "  -0.07%  "

$ws.Range("E43").Value = "  -1.79%  "

$ws.Range("D44").Value = "1.985.33"
$ws.Range("E44").Value = "  +1.19%  "

$ws.Range("D45").Value = "'0.0283"
$ws.Range("E45").Value = "  +0.44%  "

$ws.Range("E46").Value = "  +2.04%  "

$ws.Range("D47").Value = "'17.46"
$ws.Range("E47").Value = "  -1.92%  "

$ws.Range("E48").Value = "  -0.47%  "

$ws.Range("D49").Value = "2.525.96"
$ws.Range("E49").Value = "  +0.78%  "

$ws.Range("D50").Value = "'53.02"
$ws.Range("E50").Value = "  -0.95%  "

$ws.Range("E51").Value = "  -2.41%  "
